# Apply weekly fruit/vegetable price update: rows 3-8 get their
# D,I,J,K,L,M,N,P,Q values cyclically shifted (new row 3 <- old row 8,
# new row 4 <- old row 5, new row 5 <- old row 6, new row 6 <- old row 7,
# new row 7 <- old row 3, new row 8 <- old row 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target values for rows 3-8 (columns D, I, J, K, L, M, N, P, Q)
$values = @{
    3 = @{ D = 44285; I = "Primera";  J = 20; K = 25000; L = 25000; M = 25000; N = "`$/caja 18 kilos empedrada"; P = 1389; Q = 18 }
    4 = @{ D = 44280; I = "Primera";  J = 30; K = 25000; L = 25000; M = 25000; N = "`$/caja 18 kilos empedrada"; P = 1389; Q = 18 }
    5 = @{ D = 44313; I = "Primera";  J = 20; K = 15000; L = 15000; M = 15000; N = "`$/caja 15 kilos empedrada"; P = 1000; Q = 15 }
    6 = @{ D = 44313; I = "Primera";  J = 20; K = 30000; L = 30000; M = 30000; N = "`$/caja 20 kilos empedrada"; P = 1500; Q = 20 }
    7 = @{ D = 44315; I = "Especial"; J = 10; K = 30000; L = 30000; M = 30000; N = "`$/caja 20 kilos empedrada"; P = 1500; Q = 20 }
    8 = @{ D = 44315; I = "Primera";  J = 20; K = 15000; L = 15000; M = 15000; N = "`$/caja 15 kilos granel";    P = 1000; Q = 15 }
}

foreach ($row in $values.Keys) {
    $v = $values[$row]
    $ws.Range("D$row").Value = $v.D
    $ws.Range("I$row").Value = $v.I
    $ws.Range("J$row").Value = $v.J
    $ws.Range("K$row").Value = $v.K
    $ws.Range("L$row").Value = $v.L
    $ws.Range("M$row").Value = $v.M
    $ws.Range("N$row").Value = $v.N
    $ws.Range("P$row").Value = $v.P
    $ws.Range("Q$row").Value = $v.Q
}
